$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1242.2
$ws.Range("J4").Value = 1769.5
$ws.Range("L4").Value = 1769.5
$ws.Range("N4").Value = -1997.5
$ws.Range("H8").Value = 305.7857
$ws.Range("I8").Value = 252.81818
$ws.Range("J8").Value = 500
$ws.Range("K8").Value = 758.4545400000001
$ws.Range("L8").Value = 1500
$ws.Range("M8").Value = -619.4545400000001
$ws.Range("N8").Value = -1778
$ws.Range("H58").Value = 1537.1428
$ws.Range("I58").Value = 1393.3334
$ws.Range("J58").Value = 2400
$ws.Range("K58").Value = 4180.0002
$ws.Range("L58").Value = 7200
$ws.Range("M58").Value = -4030.0002
$ws.Range("N58").Value = -7500
$ws.Range("H81").Value = 37218.668
$ws.Range("J81").Value = 37218.668
$ws.Range("L81").Value = 37218.668
$ws.Range("N81").Value = -39214.668
$ws.Range("H84").Value = 37218.668
$ws.Range("J84").Value = 37218.668
$ws.Range("L84").Value = 111656.004
$ws.Range("N84").Value = -121640.004
$ws.Range("H135").Value = 107144024
$ws.Range("I135").Value = 71429200
$ws.Range("J135").Value = 142858850
$ws.Range("K135").Value = 642862800
$ws.Range("L135").Value = 1285729650
$ws.Range("M135").Value = -642860265
$ws.Range("N135").Value = -1285734720

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2478.5293
$ws.Range("I2").Value = 914.61536
$ws.Range("K2").Value = 914.61536
$ws.Range("M2").Value = -801.61536
$ws.Range("H61").Value = 11540.6
$ws.Range("I61").Value = 11849
$ws.Range("J61").Value = 10307
$ws.Range("K61").Value = 11849
$ws.Range("L61").Value = 10307
$ws.Range("M61").Value = -11637
$ws.Range("N61").Value = -10731
$ws.Range("H95").Value = 50000
$ws.Range("J95").Value = 50000
$ws.Range("L95").Value = 50000
$ws.Range("N95").Value = -55492
$ws.Range("H107").Value = 59884
$ws.Range("J107").Value = 59884
$ws.Range("L107").Value = 59884
$ws.Range("N107").Value = -67564
$ws.Range("H116").Value = 2478.5293
$ws.Range("I116").Value = 914.61536
$ws.Range("K116").Value = 914.61536
$ws.Range("M116").Value = 1379.38464
$ws.Range("H133").Value = 51189.625
$ws.Range("J133").Value = 51189.625
$ws.Range("L133").Value = 51189.625
$ws.Range("N133").Value = -56249.625
$ws.Range("H136").Value = 11540.6
$ws.Range("I136").Value = 11849
$ws.Range("J136").Value = 10307
$ws.Range("K136").Value = 35547
$ws.Range("L136").Value = 30921
$ws.Range("M136").Value = -32997
$ws.Range("N136").Value = -36021

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2478.5293
$ws.Range("I3").Value = 914.61536
$ws.Range("K3").Value = 914.61536
$ws.Range("M3").Value = -800.61536
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").ClearContents()
$ws.Range("H99").Value = 1725.7858
$ws.Range("I99").Value = 1260.091
$ws.Range("J99").Value = 3433.3333
$ws.Range("K99").Value = 1260.091
$ws.Range("L99").Value = 3433.3333
$ws.Range("M99").Value = 237.9090000000001
$ws.Range("N99").Value = -6429.3333
$ws.Range("H107").Value = 2460.9375
$ws.Range("I107").Value = 2442.4167
$ws.Range("K107").Value = 2442.4167
$ws.Range("M107").Value = -522.4167000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 476.35715
$ws.Range("I94").Value = 442
$ws.Range("J94").Value = 495.44446
$ws.Range("K94").Value = 442
$ws.Range("L94").Value = 495.44446
$ws.Range("M94").Value = 9
$ws.Range("N94").Value = -1397.44446
$ws.Range("H96").Value = 16500
$ws.Range("J96").Value = 16500
$ws.Range("L96").Value = 16500
$ws.Range("N96").Value = -21992
$ws.Range("H134").Value = 2122.1628
$ws.Range("I134").Value = 1789.4062
$ws.Range("J134").Value = 3090.182
$ws.Range("K134").Value = 5368.2186
$ws.Range("L134").Value = 9270.545999999998
$ws.Range("M134").Value = -2833.2186
$ws.Range("N134").Value = -14340.546
$ws.Range("H135").Value = 53223.75
$ws.Range("J135").Value = 53223.75
$ws.Range("L135").Value = 53223.75
$ws.Range("N135").Value = -63363.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 251.28572
$ws.Range("I11").Value = 142.33333
$ws.Range("J11").Value = 333
$ws.Range("K11").Value = 426.99999
$ws.Range("L11").Value = 999
$ws.Range("M11").Value = -286.99999
$ws.Range("N11").Value = -1279
$ws.Range("H68").Value = 2885.7415
$ws.Range("I68").Value = 864.6667
$ws.Range("J68").Value = 4032.838
$ws.Range("K68").Value = 2594.0001
$ws.Range("L68").Value = 12098.514
$ws.Range("M68").Value = -1783.0001
$ws.Range("N68").Value = -13720.514
$ws.Range("H71").Value = 2885.7415
$ws.Range("I71").Value = 864.6667
$ws.Range("J71").Value = 4032.838
$ws.Range("K71").Value = 7782.0003
$ws.Range("L71").Value = 36295.542
$ws.Range("M71").Value = -3726.0003
$ws.Range("N71").Value = -44407.542
$ws.Range("H107").Value = 920.20636
$ws.Range("I107").Value = 285.30768
$ws.Range("J107").Value = 1951.9166
$ws.Range("K107").Value = 855.92304
$ws.Range("L107").Value = 5855.7498
$ws.Range("M107").Value = 1064.07696
$ws.Range("N107").Value = -9695.7498
$ws.Range("H113").Value = 593.49207
$ws.Range("I113").Value = 572.7708
$ws.Range("J113").Value = 659.8
$ws.Range("K113").Value = 1718.3124
$ws.Range("L113").Value = 1979.4
$ws.Range("M113").Value = 451.6876
$ws.Range("N113").Value = -6319.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H31").Value = 5733.3335
$ws.Range("I31").Value = 5733.3335
$ws.Range("K31").Value = 5733.3335
$ws.Range("M31").Value = -5441.3335
$ws.Range("H37").Value = 5733.3335
$ws.Range("I37").Value = 5733.3335
$ws.Range("K37").Value = 5733.3335
$ws.Range("M37").Value = -5456.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 4000
$ws.Range("I10").Value = 4000
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 4000
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -3860
$ws.Range("N10").ClearContents()
$ws.Range("H61").Value = 32600.3
$ws.Range("I61").Value = 35900.332
$ws.Range("J61").Value = 2900
$ws.Range("K61").Value = 35900.332
$ws.Range("L61").Value = 2900
$ws.Range("M61").Value = -35698.332
$ws.Range("N61").Value = -3304
$ws.Range("H113").Value = 32600.3
$ws.Range("I113").Value = 35900.332
$ws.Range("J113").Value = 2900
$ws.Range("K113").Value = 35900.332
$ws.Range("L113").Value = 2900
$ws.Range("M113").Value = -33730.332
$ws.Range("N113").Value = -7240

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 40085.8
$ws.Range("J46").Value = 40085.8
$ws.Range("L46").Value = 40085.8
$ws.Range("N46").Value = -40547.8
$ws.Range("H69").Value = 32230.908
$ws.Range("J69").Value = 32230.908
$ws.Range("L69").Value = 32230.908
$ws.Range("N69").Value = -33728.908
$ws.Range("H72").Value = 32230.908
$ws.Range("J72").Value = 32230.908
$ws.Range("L72").Value = 96692.724
$ws.Range("N72").Value = -104180.724
$ws.Range("H126").Value = 1222.7646
$ws.Range("I126").Value = 1301.8182
$ws.Range("J126").Value = 1077.8334
$ws.Range("K126").Value = 3905.4546
$ws.Range("L126").Value = 3233.5002
$ws.Range("M126").Value = -1435.4546
$ws.Range("N126").Value = -8173.5002
$ws.Range("H134").Value = 40085.8
$ws.Range("J134").Value = 40085.8
$ws.Range("L134").Value = 120257.4
$ws.Range("N134").Value = -125327.4
$ws.Range("H136").Value = 5460.231
$ws.Range("I136").Value = 1953.1666
$ws.Range("K136").Value = 5859.4998
$ws.Range("M136").Value = -3309.4998
